$wb = $excel.ActiveWorkbook

# --- ALC row 4 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 170.73334
$ws.Range("I4").Value = 189.38461
$ws.Range("K4").Value = 189.38461
$ws.Range("M4").Value = -75.38461000000001

# --- ALC row 9 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 187.875
$ws.Range("I9").Value = 184.3077
$ws.Range("K9").Value = 184.3077
$ws.Range("M9").Value = -15.30770000000001

# --- ALC row 17 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7576302
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 7576302
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value = 22728906
$ws.Range("N17").Value = -22729242

# --- ALC row 32 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3143.3462
$ws.Range("I32").Value = 3516.8
$ws.Range("J32").Value = 2634.0908
$ws.Range("K32").Value = 3516.8
$ws.Range("L32").Value = 2634.0908
$ws.Range("M32").Value = -3190.8
$ws.Range("N32").Value = -3286.0908

# --- ALC row 70 (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5688.6665
$ws.Range("I70").Value = 10975
$ws.Range("K70").Value = 32925
$ws.Range("M70").Value = -32655

# --- ALC row 73 (hunk 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5688.6665
$ws.Range("I73").Value = 10975
$ws.Range("K73").Value = 32925
$ws.Range("M73").Value = -31989

# --- ALC row 98 (hunk 6) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 49578
$ws.Range("I98").Value = 79100.28999999999
$ws.Range("J98").Value = 20055.715
$ws.Range("K98").Value = 79100.28999999999
$ws.Range("L98").Value = 20055.715
$ws.Range("M98").Value = -77602.28999999999
$ws.Range("N98").Value = -23051.715

# --- ALC row 112 (hunk 7) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6114.4546
$ws.Range("J112").Value = 10173.167
$ws.Range("L112").Value = 30519.501
$ws.Range("N112").Value = -32735.501

# --- ALC row 122 (hunk 8) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 49578
$ws.Range("I122").Value = 79100.28999999999
$ws.Range("J122").Value = 20055.715
$ws.Range("K122").Value = 237300.87
$ws.Range("L122").Value = 60167.145
$ws.Range("M122").Value = -234850.87
$ws.Range("N122").Value = -65067.145

# --- ALC row 132 (hunk 9) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2509.4211
$ws.Range("I132").Value = 2255.453
$ws.Range("J132").Value = 5874.5
$ws.Range("K132").Value = 6766.359
$ws.Range("L132").Value = 17623.5
$ws.Range("M132").Value = -4236.359
$ws.Range("N132").Value = -22683.5

# --- ALC row 137 (hunk 10) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 288650.97
$ws.Range("I137").Value = 422321.47
$ws.Range("K137").Value = 1266964.41
$ws.Range("M137").Value = -1264414.41

# --- ARM row 61 (hunk 11) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7908.231
$ws.Range("I61").Value = 10759.2
$ws.Range("J61").Value = 6126.375
$ws.Range("K61").Value = 10759.2
$ws.Range("L61").Value = 6126.375
$ws.Range("M61").Value = -10547.2
$ws.Range("N61").Value = -6550.375

# --- ARM row 64 (hunk 12) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0

# --- ARM row 67 (hunk 13) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0

# --- ARM row 74 (hunk 14) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3580.6365
$ws.Range("I74").Value = 1857.2963
$ws.Range("K74").Value = 1857.2963
$ws.Range("M74").Value = -983.2963

# --- ARM row 77 (hunk 15) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3580.6365
$ws.Range("I77").Value = 1857.2963
$ws.Range("K77").Value = 9286.4815
$ws.Range("M77").Value = -4918.4815

# --- ARM row 98 (hunk 16) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 39355
$ws.Range("J98").Value = 39355
$ws.Range("L98").Value = 39355
$ws.Range("N98").Value = -45345

# --- ARM row 102 (hunk 17) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4075.6216
$ws.Range("J102").Value = 7849.1
$ws.Range("L102").Value = 7849.1
$ws.Range("N102").Value = -11093.1

# --- ARM row 110 (hunk 18) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 7474.8423
$ws.Range("I110").Value = 8468.134
$ws.Range("K110").Value = 8468.134
$ws.Range("M110").Value = -6423.134

# --- ARM row 132 (hunk 19) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3629.1538
$ws.Range("I132").Value = 2538.6667
$ws.Range("J132").Value = 6082.75
$ws.Range("K132").Value = 7616.000100000001
$ws.Range("L132").Value = 18248.25
$ws.Range("M132").Value = -5086.000100000001
$ws.Range("N132").Value = -23308.25

# --- ARM row 136 (hunk 20) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7908.231
$ws.Range("I136").Value = 10759.2
$ws.Range("J136").Value = 6126.375
$ws.Range("K136").Value = 32277.6
$ws.Range("L136").Value = 18379.125
$ws.Range("M136").Value = -29727.6
$ws.Range("N136").Value = -23479.125

# --- BSM row 88 (hunk 21) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 23695.428
$ws.Range("J88").Value = 25978
$ws.Range("L88").Value = 25978
$ws.Range("N88").Value = -26790

# --- BSM row 91 (hunk 22) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 23695.428
$ws.Range("J91").Value = 25978
$ws.Range("L91").Value = 25978
$ws.Range("N91").Value = -28786

# --- BSM row 134 (hunk 23) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2214.1091
$ws.Range("I134").Value = 1924.6097
$ws.Range("J134").Value = 3061.9285
$ws.Range("K134").Value = 5773.8291
$ws.Range("L134").Value = 9185.7855
$ws.Range("M134").Value = -3238.8291
$ws.Range("N134").Value = -14255.7855

# --- BSM row 138 (hunk 24) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 104996.664
$ws.Range("J138").Value = 104996.664
$ws.Range("L138").Value = 104996.664
$ws.Range("N138").Value = -115276.664

# --- CRP row 16 (hunk 25) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1640.7
$ws.Range("I16").Value = 1515.5714
$ws.Range("K16").Value = 1515.5714
$ws.Range("M16").Value = -1228.5714

# --- CRP row 28 (hunk 26) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 128650856
$ws.Range("J28").Value = 128650856
$ws.Range("L28").Value = 128650856
$ws.Range("N28").Value = -128651346

# --- CRP row 31 (hunk 27) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4909.4346
$ws.Range("I31").Value = 3950.3333
$ws.Range("K31").Value = 3950.3333
$ws.Range("M31").Value = -3655.3333

# --- CRP row 34 (hunk 28) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4909.4346
$ws.Range("I34").Value = 3950.3333
$ws.Range("K34").Value = 3950.3333
$ws.Range("M34").Value = -3748.3333

# --- CRP row 58 (hunk 29) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3255.2693
$ws.Range("I58").Value = 1910.2667
$ws.Range("K58").Value = 1910.2667
$ws.Range("M58").Value = -1707.2667

# --- CRP row 100 (hunk 30) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 78000
$ws.Range("J100").Value = 78000
$ws.Range("L100").Value = 78000
$ws.Range("N100").Value = -80164

# --- CRP row 113 (hunk 31) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1640.7
$ws.Range("I113").Value = 1515.5714
$ws.Range("K113").Value = 1515.5714
$ws.Range("M113").Value = 654.4286

# --- CRP row 132 (hunk 32) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 21576.912
$ws.Range("I132").Value = 2310.2104
$ws.Range("J132").Value = 113093.75
$ws.Range("K132").Value = 6930.6312
$ws.Range("L132").Value = 339281.25
$ws.Range("M132").Value = -4400.6312
$ws.Range("N132").Value = -344341.25

# --- CRP row 136 (hunk 33) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3255.2693
$ws.Range("I136").Value = 1910.2667
$ws.Range("K136").Value = 5730.800099999999
$ws.Range("M136").Value = -3180.800099999999

# --- CRP row 139 (hunk 34) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 46499.5
$ws.Range("J139").Value = 46499.5
$ws.Range("L139").Value = 46499.5
$ws.Range("N139").Value = -56779.5

# --- CRP row 140 (hunk 35) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 79999.5
$ws.Range("J140").Value = 79999
$ws.Range("L140").Value = 79999
$ws.Range("N140").Value = -90359

# --- CUL row 17 (hunk 36) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 423.92307
$ws.Range("I17").Value = 171.90909
$ws.Range("J17").Value = 1810
$ws.Range("K17").Value = 515.72727
$ws.Range("L17").Value = 5430
$ws.Range("M17").Value = -346.72727
$ws.Range("N17").Value = -5768

# --- CUL row 55 (hunk 37) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5187.9
$ws.Range("J55").Value = 5430.5557
$ws.Range("L55").Value = 16291.6671
$ws.Range("N55").Value = -16645.6671

# --- CUL row 108 (hunk 38) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 7824.8335
$ws.Range("I108").Value = 6237.25
$ws.Range("K108").Value = 18711.75
$ws.Range("M108").Value = -15831.75

# --- CUL row 131 (hunk 39) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 58829004
$ws.Range("J131").Value = 1783.2858
$ws.Range("L131").Value = 5349.857400000001
$ws.Range("N131").Value = -15429.8574

# --- CUL row 140 (hunk 40) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 11283.389
$ws.Range("I140").Value = 12287.5625
$ws.Range("J140").Value = 3250
$ws.Range("K140").Value = 36862.6875
$ws.Range("L140").Value = 9750
$ws.Range("M140").Value = -31682.6875
$ws.Range("N140").Value = -20110

# --- GSM row 122 (hunk 41) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 20680.812
$ws.Range("I122").Value = 18657.834
$ws.Range("K122").Value = 55973.50199999999
$ws.Range("M122").Value = -53523.50199999999

# --- GSM row 123 (hunk 42) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 18159.045
$ws.Range("J123").Value = 18159.045
$ws.Range("L123").Value = 18159.045
$ws.Range("N123").Value = -23059.045

# --- GSM row 137 (hunk 43) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 56999.5
$ws.Range("J137").Value = 56999.5
$ws.Range("L137").Value = 56999.5
$ws.Range("N137").Value = -67199.5

# --- LTW row 7 (hunk 44) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25582.736
$ws.Range("I7").Value = 36256.5
$ws.Range("K7").Value = 36256.5
$ws.Range("M7").Value = -36144.5

# --- LTW row 126 (hunk 45) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 25582.736
$ws.Range("I126").Value = 36256.5
$ws.Range("K126").Value = 108769.5
$ws.Range("M126").Value = -106299.5

# --- WVR row 98 (hunk 46) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 53000
$ws.Range("J98").Value = 53000
$ws.Range("L98").Value = 53000
$ws.Range("N98").Value = -58990

# --- WVR row 103 (hunk 47) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 39778.332
$ws.Range("J103").Value = 39778.332
$ws.Range("L103").Value = 39778.332
$ws.Range("N103").Value = -42122.332

# --- WVR row 107 (hunk 48) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 14693.739
$ws.Range("I107").Value = 1648.2354
$ws.Range("K107").Value = 4944.706200000001
$ws.Range("M107").Value = -3024.706200000001
